$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '70.699.85'
$ws.Range('E2').Value = '  +1.77%  '
$ws.Range('D3').Value = '3.556.45'
$ws.Range('E3').Value = '  +1.64%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '585.04'
$ws.Range('E5').Value = '  +0.23%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '187.99'
$ws.Range('E6').Value = '  +2.60%  '
$ws.Range('D7').Value = '3.548.44'
$ws.Range('E7').Value = '  +1.69%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.620'
$ws.Range('E8').Value = '  +1.61%  '
$ws.Range('E9').Value = '  -0.09%  '
$ws.Range('E10').Value = '  +6.91%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.642'
$ws.Range('E11').Value = '  -0.26%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '53.81'
$ws.Range('E12').Value = '  +0.22%  '
$ws.Range('E13').Value = '  +0.94%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '9.38'
$ws.Range('E14').Value = '  -0.38%  '
$ws.Range('D15').Value = '4.123.92'
$ws.Range('E15').Value = '  +1.67%  '
$ws.Range('B16').Value = 'WrappedBTC'
$ws.Range('C16').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D16').Value = '70.721.32'
$ws.Range('E16').Value = '  +1.87%  '
$ws.Range('B17').Value = 'WrappedEther'
$ws.Range('C17').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D17').Value = '3.542.48'
$ws.Range('E17').Value = '  +1.23%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '12.65'
$ws.Range('E18').Value = '  +2.84%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '18.91'
$ws.Range('E19').Value = '  -1.79%  '
$ws.Range('B20').Value = 'TRON'
$ws.Range('C20').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '0.120'
$ws.Range('E20').Value = '  +0.75%  '
$ws.Range('B21').Value = 'BitcoinCash'
$ws.Range('C21').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '562.37'
$ws.Range('E21').Value = '  +3.58%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.991'
$ws.Range('E22').Value = '  -1.72%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '17.77'
$ws.Range('E23').Value = '  -3.93%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '4.59'
$ws.Range('E24').Value = '  +1.33%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '93.69'
$ws.Range('E26').Value = '  -1.54%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '11.03'
$ws.Range('E27').Value = '  -0.71%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '2.90'
$ws.Range('E28').Value = '  -2.45%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '9.26'
$ws.Range('E29').Value = '  +1.67%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '32.26'
$ws.Range('E30').Value = '  +1.49%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '7.02'
$ws.Range('E31').Value = '  -3.08%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '12.15'
$ws.Range('E32').Value = '  -2.89%  '
$ws.Range('E33').Value = '  +26.50%  '
$ws.Range('E34').Value = '  +1.44%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '63.06'
$ws.Range('E35').Value = '  -1.11%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '3.23'
$ws.Range('E36').Value = '  +5.15%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '528.36'
$ws.Range('E37').Value = '  +0.33%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.404'
$ws.Range('E38').Value = '  -0.67%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '37.86'
$ws.Range('E39').Value = '  -0.31%  '
$ws.Range('B40').Value = 'Dai'
$ws.Range('C40').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.999'
$ws.Range('E40').Value = '  +0.04%  '
$ws.Range('B41').Value = 'Maker'
$ws.Range('C41').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D41').Value = '3.636.17'
$ws.Range('E41').Value = '  +9.10%  '
$ws.Range('D42').Value = '0.0₃0784'
$ws.Range('E42').Value = '  +2.91%  '
$ws.Range('E43').Value = '  +4.65%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.138'
$ws.Range('E44').Value = '  +2.86%  '
$ws.Range('E45').Value = '  +3.50%  '
$ws.Range('E46').Value = '  -0.52%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '2.91'
$ws.Range('E47').Value = '  -1.98%  '
$ws.Range('E48').Value = '  +2.61%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '9.16'
$ws.Range('E49').Value = '  +2.17%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '1.44'
$ws.Range('E51').Value = '  +8.75%  '
